# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# The per-start strikeout counts below were recomputed/re-sourced and written
# back into the "K" column (column G) of the single data sheet, row by row
# (row 1 is the header row; data starts at row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 1
    6 = 2
    7 = 0
    8 = 0
    9 = 0
    10 = 1
    11 = 1
    12 = 3
    13 = 1
    14 = 0
    15 = 2
    16 = 3
    17 = 3
    18 = 2
    19 = 1
    20 = 1
    21 = 1
    22 = 2
    23 = 3
    24 = 1
    25 = 2
    26 = 1
    27 = 1
    28 = 2
    29 = 1
    30 = 1
    31 = 1
    32 = 3
    33 = 1
    34 = 2
    35 = 2
    36 = 3
    37 = 2
    38 = 0
    39 = 0
    40 = 0
    41 = 3
    42 = 2
    43 = 1
    44 = 2
    45 = 1
    46 = 3
    47 = 0
    48 = 2
    49 = 1
    50 = 0
    51 = 1
    52 = 1
    53 = 0
    54 = 1
    55 = 1
    56 = 2
    57 = 1
    58 = 2
    59 = 3
    60 = 2
    61 = 1
    62 = 2
    63 = 2
    64 = 2
    67 = 1
    68 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Host "Updated K column (G) for" $kValues.Count "rows"
